$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last "Total" row (row 36) only keeps its Product/Total
# label and the Subtotal value once it is no longer the most recent
# order; the placeholder blanks in A/B/D/E/F are dropped.
$ws.Cells.Item(36, 1).ClearContents()
$ws.Cells.Item(36, 2).ClearContents()
$ws.Cells.Item(36, 4).ClearContents()
$ws.Cells.Item(36, 5).ClearContents()
$ws.Cells.Item(36, 6).ClearContents()

# New order row (row 37)
$ws.Cells.Item(37, 1).Value = "2025-03-01 01:57:52"
$ws.Cells.Item(37, 2).Value = "benichi"
$ws.Cells.Item(37, 3).Value = "Sencilla"
$ws.Cells.Item(37, 4).Value = "Salchipapas"
$ws.Cells.Item(37, 5).Value = 2
$ws.Cells.Item(37, 6).Value = 10000
$ws.Cells.Item(37, 7).Value = 20000

# Total row for the new order (row 38)
$ws.Cells.Item(38, 3).Value = "Total"
$ws.Cells.Item(38, 7).Value = 20000
